$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '50.746.65'
$ws.Range("E2").Value = '  -1.20%  '
$ws.Range("D3").Value = '2.920.72'
$ws.Range("E3").Value = '  -1.80%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'374.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.12%  '
$ws.Range("D6").Value = "'99.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.29%  '
$ws.Range("E7").Value = '  -1.16%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'0.576"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("D10").Value = "'35.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.75%  '
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").Value = "'0.0845"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").Value = '3.377.77'
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("D14").Value = "'17.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("D15").Value = "'7.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("D16").Value = "'11.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +62.73%  '
$ws.Range("D17").Value = '2.927.26'
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").Value = "'0.993"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = '50.702.32'
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("D20").Value = "'3.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.43%  '
$ws.Range("D21").Value = "'12.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.67%  '
$ws.Range("E22").Value = '  -1.65%  '
$ws.Range("D23").Value = "'69.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("D24").Value = "'265.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("E25").Value = '  +8.37%  '
$ws.Range("D26").Value = "'7.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.00%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = "'7.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.04%  '
$ws.Range("D29").Value = "'25.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("E30").Value = '  -3.17%  '
$ws.Range("E31").Value = '  -4.77%  '
$ws.Range("D32").Value = "'9.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").Value = "'50.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("D35").Value = "'33.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.11%  '
$ws.Range("E36").Value = '  -3.91%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("D40").Value = "'16.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.64%  '
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("E42").Value = '  -6.41%  '
$ws.Range("D43").Value = "'119.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.45%  '
$ws.Range("D44").Value = "'20.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.11%  '
$ws.Range("E45").Value = '  -1.99%  '
$ws.Range("D46").Value = "'3.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.64%  '
$ws.Range("D48").Value = '1.989.57'
$ws.Range("E48").Value = '  -1.70%  '
$ws.Range("D49").Value = "'0.258"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.49%  '
$ws.Range("D50").Value = "'0.0313"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.26%  '
$ws.Range("D51").Value = "'5.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.97%  '
